$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saisie")
$ws.Range("A1").Value = "TEST"
